# Generate Report for Handoff
# Updates the "b.md" row (row 3) across all three sheets (Overview, zh-cn, de-de)
# to reflect that a new handoff has been generated for b.md.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Overview"
#   Columns: A=File Name, B=Path And Name, C=Extension, D=Publish URL,
#            E=zh-cn, F=de-de, G=Latest HO Xliff Generate Date
#   Row 3 = b.md
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(3, 5).Value = "Ready for handoff"          # E3 zh-cn
$wsOverview.Cells.Item(3, 6).Value = "Ready for handoff"          # F3 de-de
$wsOverview.Cells.Item(3, 7).Value = "2016-09-03 08:41:09"        # G3 Latest HO Xliff Generate Date

# ---------------------------------------------------------------
# Sheet 2: "zh-cn"
#   Columns: A=Source File Name, B=File Extension, C=Status, D=Source Path,
#            E=Priority, F=Content Duplicate, G=Latest Handoff File,
#            H=Latest Handoff Datetime, I=Latest Target File, J=Latest Handback File,
#            K=Latest Handback DateTime, L=Reference Tokens, M=To be localized,
#            N=Dependency From, O=Has metadata, P=Error Detail
#   Row 3 = b.md
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Item(3, 3).Value = "Ready for handoff"               # C3 Status

# Content Duplicate (F3) must remain stored as TEXT "False" (matching the rest
# of the sheet), not the Boolean FALSE that a bare assignment would produce.
# A leading apostrophe forces text entry; resetting the style afterwards drops
# the "quote prefix" formatting that the apostrophe trick applies.
$wsZhCn.Cells.Item(3, 6).Value = "'False"                          # F3 Content Duplicate
$wsZhCn.Cells.Item(3, 6).Style = "Normal"

$wsZhCn.Cells.Item(3, 7).Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"  # G3 Latest Handoff File
$wsZhCn.Cells.Item(3, 8).Value = "2016-09-03 08:40:59"             # H3 Latest Handoff Datetime
$wsZhCn.Cells.Item(3, 16).Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dbd92ffafc8011c77d819337ed6b2c91e3bbac65/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ee1ae0d480e387fb3e9286b27d4e6dec782c498f/e2e/b.md."  # P3 Error Detail

# Widen the Error Detail column (P) to match the new, longer text.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.15

# ---------------------------------------------------------------
# Sheet 3: "de-de"
#   Same column layout as zh-cn.
#   Row 3 = b.md
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Item(3, 3).Value = "Ready for handoff"               # C3 Status

$wsDeDe.Cells.Item(3, 6).Value = "'False"                          # F3 Content Duplicate
$wsDeDe.Cells.Item(3, 6).Style = "Normal"

$wsDeDe.Cells.Item(3, 7).Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"  # G3 Latest Handoff File
$wsDeDe.Cells.Item(3, 8).Value = "2016-09-03 08:41:09"             # H3 Latest Handoff Datetime
$wsDeDe.Cells.Item(3, 16).Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dbd92ffafc8011c77d819337ed6b2c91e3bbac65/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ee1ae0d480e387fb3e9286b27d4e6dec782c498f/e2e/b.md."  # P3 Error Detail

# Widen the Error Detail column (P) to match the new, longer text.
$wsDeDe.Columns.Item(16).ColumnWidth = 39.15
